$d = $word.ActiveDocument

# Replace the run content of an existing paragraph (in place) with one or
# more runs of text, preserving the paragraph's own <w:p> attributes
# (w14:paraId, rsidR, etc.) and its <w:pPr> (pStyle/numPr at the given
# ilvl). Runs whose text has leading/trailing whitespace get
# xml:space="preserve" so Word doesn't trim it on load.
function Set-ParagraphRuns($paraIndex, $ilvl, $runs) {
    $p = $d.Paragraphs.Item($paraIndex)
    # Stop one character short of the paragraph mark so InsertXML replaces
    # only the paragraph's content, not the mark itself (replacing through
    # the mark on the very last paragraph of the body leaves a stray empty
    # paragraph behind).
    $r = $d.Range($p.Range.Start, $p.Range.End - 1)

    $runXml = ""
    foreach ($run in $runs) {
        if ($run -match '^\s' -or $run -match '\s$') {
            $runXml += "<w:r><w:t xml:space=`"preserve`">" + $run + "</w:t></w:r>"
        } else {
            $runXml += "<w:r><w:t>" + $run + "</w:t></w:r>"
        }
    }

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="' + $ilvl + '"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
        $runXml + '</w:p>' +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($xml) | Out-Null
}

# Insert a brand-new ListParagraph-styled paragraph right after the
# paragraph at paraIndex, at the given list level (0 or 1).
function Insert-ParagraphAfter($paraIndex, $ilvl, $text) {
    $p = $d.Paragraphs.Item($paraIndex)
    $insertPos = $p.Range.End - 1
    $r = $d.Range($insertPos, $insertPos)

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="' + $ilvl + '"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
        '<w:r><w:t>' + $text + '</w:t></w:r></w:p>' +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($xml) | Out-Null
}

# 1) "Add computer logic when it's the last space" (ilvl 1, paragraph 4)
#    -> two runs: the "how many are left" sentence + " (more coins at end)"
Set-ParagraphRuns 4 "1" @(
    "Add “how many are left” to logic so comp uses coins",
    " (more coins at end)"
)

# 2) "Add "how many are left" to logic so comp uses coins" (ilvl 1, paragraph 5)
#    -> replaced by the "destroying your tiles" sentence
Set-ParagraphRuns 5 "1" @(
    "Add “… destroying your tiles on __ and ___.”"
)

# New paragraph "PHP" (ilvl 0) inserted right after it
Insert-ParagraphAfter 5 "0" "PHP"

# 3) "Compute to use bomb better (like when it would help him)" is now
#    paragraph 7 (shifted down by the PHP insertion above)
#    -> replaced by the "Save user information" sentence
Set-ParagraphRuns 7 "1" @(
    "Save user information to display at the end (with ability to delete)"
)

# New paragraph "Count up score at end" (ilvl 0) inserted right after it
Insert-ParagraphAfter 7 "0" "Count up score at end"
